$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1934.919849777034,
    2988.576002886201,
    4166.878890431921,
    4696.545507286738,
    5102.612874897974,
    5936.571654466956,
    6231.96760524919,
    6409.291201575541,
    6750.207609667235,
    6953.238624305999,
    7072.01671505661,
    7408.09351793669,
    7560.113295091199,
    7672.016297393805,
    7873.389020378461,
    8008.031782703364,
    8095.19391187989,
    8167.362054571519,
    8341.539911358934,
    8390.094114595906,
    8453.696374182065,
    8457.615102224545,
    8460.208531760383,
    8467.217522703546,
    8460.199749823316,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597,
    8460.308567009597
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
